$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename column headers in row 1:
#    "<field>_old" -> "<field>_FV2410"   (columns A-J, 1-10)
#    "diff"        stays "diff"          (column K, 11)
#    "<field>_new" -> "<field>_FV2504"   (columns L-U, 12-21)
$oldSuffixHeaders = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)
$newSuffixHeaders = @(
    "Segmentname_FV2504",
    "Segmentgruppe_FV2504",
    "Segment_FV2504",
    "Datenelement_FV2504",
    "Segment ID_FV2504",
    "Code_FV2504",
    "Qualifier_FV2504",
    "Beschreibung_FV2504",
    "Bedingungsausdruck_FV2504",
    "Bedingung_FV2504"
)

for ($i = 0; $i -lt $oldSuffixHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $oldSuffixHeaders[$i]
}
$ws.Cells.Item(1, 11).Value = "diff"
for ($i = 0; $i -lt $newSuffixHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $newSuffixHeaders[$i]
}

# 2. Turn the used range into an Excel Table ("Table1") over A1:U63,
#    picking up the freshly renamed header row for the table column names.
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U63"), [System.Reflection.Missing]::Value, 1)
$tbl.Name = "Table1"

# 3. Freeze the header row (split below row 1, top-left of the scrolling
#    pane at A2).
$ws.Range("A2").Select() | Out-Null
($excel.ActiveWindow.FreezePanes = $true) | Out-Null
